$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.070.71"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "2.432.31"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.54"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.22"
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "2.418.73"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D16").Value = "2.820.14"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "61.036.39"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "2.427.08"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.14"
$ws.Range("E21").Value = "  -2.50%  "
$ws.Range("E22").Value = "  -2.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.11"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.88"
$ws.Range("E25").Value = "  -5.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.91"
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.73"
$ws.Range("E27").Value = "  -7.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "574.25"
$ws.Range("E28").Value = "  -7.73%  "
$ws.Range("D29").Value = "2.542.44"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("D30").Value = "0.0₃0910"
$ws.Range("E30").Value = "  -5.01%  "
$ws.Range("E31").Value = "  -2.15%  "
$ws.Range("E32").Value = "  -6.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  -2.15%  "
$ws.Range("E34").Value = "  -6.70%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  -6.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.368"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "148.81"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.07"
$ws.Range("E41").Value = "  -4.18%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.71"
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.66"
$ws.Range("E44").Value = "  -5.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("E45").Value = "  -6.40%  "
$ws.Range("D46").Value = "0.0₆0283"
$ws.Range("E46").Value = "  +18.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.79"
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.592"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("E50").Value = "  -4.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.32"
$ws.Range("E51").Value = "  -1.40%  "
